# The commit swaps which sheet/tab is active & selected when the workbook is
# opened: previously "ONBRA" was the active/selected tab, now "PROXIES" is
# the active/selected tab (scrolled so row 2 is at the top). The other
# sheets' own selections are left untouched.

$wb = $excel.ActiveWorkbook

$proxies = $wb.Worksheets.Item("PROXIES")

# Make PROXIES the active sheet (this also clears "tabSelected" from
# whichever sheet was previously active, i.e. ONBRA, and updates the
# workbook-level active-tab pointer).
$proxies.Activate()

# Scroll PROXIES so that row 2 is the top visible row (keeps the existing
# selection, L9, untouched).
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
